$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E13:E21 previously held "Yes" in a bold font while E2:E12 held "No" in
# a plain (non-bold) font. Make the whole column uniform like the former
# "No" cells (non-bold) and set every cell in E2:E21 to "Yes".
$ws.Range("E13:E21").Font.Bold = $false

$ws.Range("E2:E21").Value = "Yes"

$ws.Range("E2:E21").Select()
